$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsBoundaries = $wb.Worksheets.Item("Boundaries and methane sources")

# --- About sheet ---
$wsAbout.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

$wsAbout.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for Karvina Coal Mines, Czech Republic, M0449, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)" + [char]34

# --- Boundaries and methane sources sheet ---
for ($r = 2; $r -le 9; $r++) {
    $cell = $wsBoundaries.Range("S$r")
    $current = $cell.Value()
    if ($current -ne $null -and $current.Contains($oldStamp)) {
        $cell.Value = $current.Replace($oldStamp, $newStamp)
    }
}
